$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for new column D - copy formatting (bold, border, alignment) from C1
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "%_cumulative"

# Cumulative percentage values for column D (rows 2-14)
$values = @(
    36.51,
    67.40000000000001,
    83.71000000000001,
    91.79000000000001,
    95.78,
    97.91,
    99.03999999999999,
    99.59999999999999,
    99.80999999999999,
    99.91999999999999,
    99.97999999999999,
    99.98999999999999,
    99.98999999999999
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}
